$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (ALC)
$ws.Range("H19").Value = 1390
$ws.Range("I19").Value = 1258.2
$ws.Range("J19").Value = 1499.8334
$ws.Range("K19").Value = 1258.2
$ws.Range("L19").Value = 1499.8334
$ws.Range("M19").Value = -1083.2
$ws.Range("N19").Value = -1849.8334

# Row 48 (ALC)
$ws.Range("H48").Value = 3000
$ws.Range("J48").Value = 3000
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9584

# Row 56 (ALC)
$ws.Range("H56").Value = 3000
$ws.Range("J56").Value = 3000
$ws.Range("L56").Value = 9000
$ws.Range("N56").Value = -10068

# Row 86 (ALC)
$ws.Range("H86").Value = 2268.25
$ws.Range("I86").Value = 2236.6428
$ws.Range("J86").Value = 2299.8572
$ws.Range("K86").Value = 2236.6428
$ws.Range("L86").Value = 2299.8572
$ws.Range("M86").Value = -1113.6428
$ws.Range("N86").Value = -4545.8572

# Row 89 (ALC)
$ws.Range("H89").Value = 2268.25
$ws.Range("I89").Value = 2236.6428
$ws.Range("J89").Value = 2299.8572
$ws.Range("K89").Value = 11183.214
$ws.Range("L89").Value = 11499.286
$ws.Range("M89").Value = -5567.214
$ws.Range("N89").Value = -22731.286

# Row 113 (ALC)
$ws.Range("H113").Value = 64464.562
$ws.Range("I113").Value = 113092.664
$ws.Range("J113").Value = 1942.7142
$ws.Range("K113").Value = 113092.664
$ws.Range("L113").Value = 1942.7142
$ws.Range("M113").Value = -109838.664
$ws.Range("N113").Value = -8450.7142

# Row 129 (ALC)
$ws.Range("H129").Value = 3670.8333
$ws.Range("J129").Value = 933.1818
$ws.Range("L129").Value = 2799.5454
$ws.Range("N129").Value = -12799.5454

# Row 131 (ALC)
$ws.Range("H131").Value = 4010.5686
$ws.Range("I131").Value = 942.9
$ws.Range("J131").Value = 4758.7803
$ws.Range("K131").Value = 2828.7
$ws.Range("L131").Value = 14276.3409
$ws.Range("M131").Value = 2211.3
$ws.Range("N131").Value = -24356.3409

# Row 132 (ALC)
$ws.Range("H132").Value = 5006494
$ws.Range("I132").Value = 5820644
$ws.Range("J132").Value = 5285.5713
$ws.Range("K132").Value = 17461932
$ws.Range("L132").Value = 15856.7139
$ws.Range("M132").Value = -17459402
$ws.Range("N132").Value = -20916.7139

# Row 137 (ALC)
$ws.Range("H137").Value = 1804.5416
$ws.Range("I137").Value = 1600.5385
$ws.Range("J137").Value = 2045.6364
$ws.Range("K137").Value = 4801.6155
$ws.Range("L137").Value = 6136.9092
$ws.Range("M137").Value = -2251.6155
$ws.Range("N137").Value = -11236.9092

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 8729.392
$ws.Range("I32").Value = 7385.9775
$ws.Range("J32").Value = 26001.857
$ws.Range("K32").Value = 7385.9775
$ws.Range("L32").Value = 26001.857
$ws.Range("M32").Value = -7098.9775
$ws.Range("N32").Value = -26575.857

# Row 45 (ARM)
$ws.Range("H45").Value = 51513.05
$ws.Range("I45").Value = 78055.69500000001
$ws.Range("J45").Value = 2219.5715
$ws.Range("K45").Value = 78055.69500000001
$ws.Range("L45").Value = 2219.5715
$ws.Range("M45").Value = -77678.69500000001
$ws.Range("N45").Value = -2973.5715

# Row 74 (ARM)
$ws.Range("H74").Value = 1111.2354
$ws.Range("I74").Value = 1102.2609
$ws.Range("K74").Value = 1102.2609
$ws.Range("M74").Value = -228.2609

# Row 77 (ARM)
$ws.Range("H77").Value = 1111.2354
$ws.Range("I77").Value = 1102.2609
$ws.Range("K77").Value = 5511.3045
$ws.Range("M77").Value = -1143.3045

# Row 97 (ARM)
$ws.Range("H97").Value = 64518.25
$ws.Range("I97").Value = 144298.58
$ws.Range("J97").Value = 2466.889
$ws.Range("K97").Value = 144298.58
$ws.Range("L97").Value = 2466.889
$ws.Range("M97").Value = -143802.58
$ws.Range("N97").Value = -3458.889

# Row 110 (ARM)
$ws.Range("H110").Value = 50106196
$ws.Range("I110").Value = 52743304
$ws.Range("J110").Value = 1200
$ws.Range("K110").Value = 52743304
$ws.Range("L110").Value = 1200
$ws.Range("M110").Value = -52741259
$ws.Range("N110").Value = -5290

# Row 122 (ARM)
$ws.Range("H122").Value = 2417
$ws.Range("I122").Value = 2019.6875
$ws.Range("K122").Value = 6059.0625
$ws.Range("M122").Value = -3609.0625

# Row 132 (ARM)
$ws.Range("H132").Value = 17000.514
$ws.Range("I132").Value = 29340.904
$ws.Range("J132").Value = 2603.389
$ws.Range("K132").Value = 88022.712
$ws.Range("L132").Value = 7810.167
$ws.Range("M132").Value = -85492.712
$ws.Range("N132").Value = -12870.167

$ws = $wb.Worksheets.Item("BSM")
# Row 8 (BSM)
$ws.Range("H8").Value = 27200
$ws.Range("I8").Value = 2000
$ws.Range("J8").Value = 39800
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 39800
$ws.Range("M8").Value = -1860
$ws.Range("N8").Value = -40080

# Row 33 (BSM)
$ws.Range("H33").Value = 5600
$ws.Range("I33").Value = 4000
$ws.Range("K33").Value = 4000
$ws.Range("M33").Value = -3664

# Row 94 (BSM)
$ws.Range("H94").Value = 628.1429000000001
$ws.Range("I94").Value = 532.8333
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 532.8333
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -81.83330000000001
$ws.Range("N94").Value = -2102

# Row 99 (BSM)
$ws.Range("H99").Value = 2225
$ws.Range("I99").Value = 3580
$ws.Range("J99").Value = 2074.4443
$ws.Range("K99").Value = 3580
$ws.Range("L99").Value = 2074.4443
$ws.Range("M99").Value = -2082
$ws.Range("N99").Value = -5070.4443

# Row 107 (BSM)
$ws.Range("H107").Value = 41686860
$ws.Range("I107").Value = 52656724
$ws.Range("J107").Value = 1368
$ws.Range("K107").Value = 52656724
$ws.Range("L107").Value = 1368
$ws.Range("M107").Value = -52654804
$ws.Range("N107").Value = -5208

# Row 134 (BSM)
$ws.Range("H134").Value = 3388.825
$ws.Range("I134").Value = 3389.4333
$ws.Range("J134").Value = 3387
$ws.Range("K134").Value = 10168.2999
$ws.Range("L134").Value = 10161
$ws.Range("M134").Value = -7633.2999
$ws.Range("N134").Value = -15231

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 29973.725
$ws.Range("I31").Value = 841.2593000000001
$ws.Range("K31").Value = 841.2593000000001
$ws.Range("M31").Value = -546.2593000000001

# Row 34 (CRP)
$ws.Range("H34").Value = 29973.725
$ws.Range("I34").Value = 841.2593000000001
$ws.Range("K34").Value = 841.2593000000001
$ws.Range("M34").Value = -639.2593000000001

# Row 105 (CRP)
$ws.Range("H105").Value = 1104.75
$ws.Range("I105").Value = 1091.5
$ws.Range("J105").Value = 1144.5
$ws.Range("K105").Value = 1091.5
$ws.Range("L105").Value = 1144.5
$ws.Range("M105").Value = 655.5
$ws.Range("N105").Value = -4638.5

# Row 132 (CRP)
$ws.Range("H132").Value = 2422.6
$ws.Range("I132").Value = 2305.5938
$ws.Range("K132").Value = 6916.7814
$ws.Range("M132").Value = -4386.7814

# Row 134 (CRP)
$ws.Range("H134").Value = 2132.0908
$ws.Range("I134").Value = 970
$ws.Range("J134").Value = 3100.5
$ws.Range("K134").Value = 2910
$ws.Range("L134").Value = 9301.5
$ws.Range("M134").Value = -375
$ws.Range("N134").Value = -14371.5

$ws = $wb.Worksheets.Item("CUL")
# Row 37 (CUL)
$ws.Range("H37").Value = 855338.5
$ws.Range("J37").Value = 855338.5
$ws.Range("L37").Value = 2566015.5
$ws.Range("N37").Value = -2566239.5

# Row 92 (CUL)
$ws.Range("H92").Value = 1000
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 3000
$ws.Range("N92").Value = -5496

# Row 96 (CUL)
$ws.Range("H96").Value = 21833.334
$ws.Range("J96").Value = 21833.334
$ws.Range("L96").Value = 65500.00199999999
$ws.Range("N96").Value = -69618.00199999999

# Row 101 (CUL)
$ws.Range("H101").Value = 4081.4546
$ws.Range("J101").Value = 4081.4546
$ws.Range("L101").Value = 12244.3638
$ws.Range("N101").Value = -17112.3638

# Row 102 (CUL)
$ws.Range("H102").Value = 4874
$ws.Range("I102").Value = 4500
$ws.Range("J102").Value = 4998.6665
$ws.Range("K102").Value = 13500
$ws.Range("L102").Value = 14995.9995
$ws.Range("M102").Value = -11066
$ws.Range("N102").Value = -19863.9995

# Row 110 (CUL)
$ws.Range("H110").Value = 1580
$ws.Range("I110").Value = 1580
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 4740
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -650
$ws.Range("N110").ClearContents()

# Row 111 (CUL)
$ws.Range("H111").Value = 1000
$ws.Range("I111").Value = 1000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3000
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 67
$ws.Range("N111").ClearContents()

# Row 131 (CUL)
$ws.Range("H131").Value = 1037.06
$ws.Range("I131").Value = 428.22726
$ws.Range("J131").Value = 1208.7821
$ws.Range("K131").Value = 1284.68178
$ws.Range("L131").Value = 3626.3463
$ws.Range("M131").Value = 3755.31822
$ws.Range("N131").Value = -13706.3463

# Row 137 (CUL)
$ws.Range("H137").Value = 5057872
$ws.Range("I137").Value = 93124.55
$ws.Range("J137").Value = 11125896
$ws.Range("K137").Value = 279373.65
$ws.Range("L137").Value = 33377688
$ws.Range("M137").Value = -274273.65
$ws.Range("N137").Value = -33387888

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Range("H122").Value = 3350.8
$ws.Range("I122").Value = 2562.5
$ws.Range("K122").Value = 7687.5
$ws.Range("M122").Value = -5237.5

# Row 126 (GSM)
$ws.Range("H126").Value = 4528409.5
$ws.Range("I126").Value = 3836.3333
$ws.Range("J126").Value = 14708699
$ws.Range("K126").Value = 11508.9999
$ws.Range("L126").Value = 44126097
$ws.Range("M126").Value = -9038.999899999999
$ws.Range("N126").Value = -44131037

# Row 132 (GSM)
$ws.Range("H132").Value = 4158.579
$ws.Range("I132").Value = 3310.8
$ws.Range("J132").Value = 5100.5557
$ws.Range("K132").Value = 9932.400000000001
$ws.Range("L132").Value = 15301.6671
$ws.Range("M132").Value = -7402.400000000001
$ws.Range("N132").Value = -20361.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("H16").Value = 689787.5
$ws.Range("I16").Value = 144485.86
$ws.Range("J16").Value = 1113911
$ws.Range("K16").Value = 144485.86
$ws.Range("L16").Value = 1113911
$ws.Range("M16").Value = -144315.86
$ws.Range("N16").Value = -1114251

# Row 40 (LTW)
$ws.Range("H40").Value = 50592.715
$ws.Range("I40").Value = 114132.445
$ws.Range("J40").Value = 2937.9167
$ws.Range("K40").Value = 114132.445
$ws.Range("L40").Value = 2937.9167
$ws.Range("M40").Value = -113996.445
$ws.Range("N40").Value = -3209.9167

# Row 132 (LTW)
$ws.Range("H132").Value = 4069.5417
$ws.Range("I132").Value = 4404.4707
$ws.Range("J132").Value = 3256.1428
$ws.Range("K132").Value = 13213.4121
$ws.Range("L132").Value = 9768.428400000001
$ws.Range("M132").Value = -10683.4121
$ws.Range("N132").Value = -14828.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 7 (WVR)
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
